$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 283.5
$ws.Range("I2").Value = 185.7
$ws.Range("K2").Value = 185.7
$ws.Range("M2").Value = -72.69999999999999

# ALC row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 6510
$ws.Range("I21").Value = 6883.6665
$ws.Range("K21").Value = 6883.6665
$ws.Range("M21").Value = -6415.6665

# ALC row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 6510
$ws.Range("I23").Value = 6883.6665
$ws.Range("K23").Value = 6883.6665
$ws.Range("M23").Value = -6649.6665

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1996.2727
$ws.Range("I38").Value = 324.2857
$ws.Range("J38").Value = 4922.25
$ws.Range("K38").Value = 972.8571000000001
$ws.Range("L38").Value = 14766.75
$ws.Range("M38").Value = -600.8571000000001
$ws.Range("N38").Value = -15510.75

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1672.5555
$ws.Range("I58").Value = 138.25
$ws.Range("J58").Value = 2900
$ws.Range("K58").Value = 414.75
$ws.Range("L58").Value = 8700
$ws.Range("M58").Value = -264.75
$ws.Range("N58").Value = -9000

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7362.125
$ws.Range("I64").Value = 5974.25
$ws.Range("J64").Value = 8750
$ws.Range("K64").Value = 5974.25
$ws.Range("L64").Value = 8750
$ws.Range("M64").Value = -5726.25
$ws.Range("N64").Value = -9246

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 7362.125
$ws.Range("I67").Value = 5974.25
$ws.Range("J67").Value = 8750
$ws.Range("K67").Value = 5974.25
$ws.Range("L67").Value = 8750
$ws.Range("M67").Value = -5116.25
$ws.Range("N67").Value = -10466

# ALC row 68
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = ""
$ws.Range("N68").Value = 0

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 6391.5
$ws.Range("I69").Value = 5013
$ws.Range("K69").Value = 15039
$ws.Range("M69").Value = -14165

# ALC row 71
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = ""
$ws.Range("N71").Value = 0

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 6391.5
$ws.Range("I72").Value = 5013
$ws.Range("K72").Value = 45117
$ws.Range("M72").Value = -40749

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1849.7646
$ws.Range("I100").Value = 715.1
$ws.Range("J100").Value = 3470.7144
$ws.Range("K100").Value = 715.1
$ws.Range("L100").Value = 3470.7144
$ws.Range("M100").Value = -174.1
$ws.Range("N100").Value = -4552.7144

# ALC row 104
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H104").Value = 1021.5
$ws.Range("I104").Value = 1025.8
$ws.Range("J104").Value = 1000
$ws.Range("K104").Value = 3077.4
$ws.Range("L104").Value = 3000
$ws.Range("M104").Value = -1330.4
$ws.Range("N104").Value = -6494

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1970
$ws.Range("I113").Value = 1858
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 1858
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 1396
$ws.Range("N113").Value = -8758

# ARM row 25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 19750
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = ""

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2585.625
$ws.Range("I61").Value = 2585.625
$ws.Range("K61").Value = 2585.625
$ws.Range("M61").Value = -2373.625

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3049.5
$ws.Range("I132").Value = 3049.5
$ws.Range("K132").Value = 9148.5
$ws.Range("M132").Value = -6618.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2585.625
$ws.Range("I136").Value = 2585.625
$ws.Range("K136").Value = 7756.875
$ws.Range("M136").Value = -5206.875

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3751.7058
$ws.Range("I86").Value = 2124.7273
$ws.Range("K86").Value = 2124.7273
$ws.Range("M86").Value = -1001.7273

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3751.7058
$ws.Range("I89").Value = 2124.7273
$ws.Range("K89").Value = 10623.6365
$ws.Range("M89").Value = -5007.636500000001

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1174.8125
$ws.Range("I22").Value = 522.6667
$ws.Range("K22").Value = 522.6667
$ws.Range("M22").Value = -172.6667

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2699.1
$ws.Range("I99").Value = 2141.4285
$ws.Range("K99").Value = 2141.4285
$ws.Range("M99").Value = -643.4285

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1391.4615
$ws.Range("I105").Value = 1432.4166
$ws.Range("K105").Value = 1432.4166
$ws.Range("M105").Value = 314.5834

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2699.1
$ws.Range("I126").Value = 2141.4285
$ws.Range("K126").Value = 6424.2855
$ws.Range("M126").Value = -3954.2855

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 68295.2
$ws.Range("J141").Value = 68295.2
$ws.Range("L141").Value = 68295.2
$ws.Range("N141").Value = -78655.2

# CUL row 25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 6
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = ""

# CUL row 30
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 6
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = ""

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 8257.143
$ws.Range("J39").Value = 9500
$ws.Range("L39").Value = 28500
$ws.Range("N39").Value = -29088

# CUL row 54
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2743.75
$ws.Range("J54").Value = 2733.3333
$ws.Range("L54").Value = 8199.999899999999
$ws.Range("N54").Value = -9317.999899999999

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 6160.2
$ws.Range("J55").Value = 7312.25
$ws.Range("L55").Value = 21936.75
$ws.Range("N55").Value = -22290.75

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2363.3333
$ws.Range("I136").Value = 636
$ws.Range("J136").Value = 11000
$ws.Range("K136").Value = 1908
$ws.Range("L136").Value = 33000
$ws.Range("M136").Value = 3192
$ws.Range("N136").Value = -43200

# GSM row 24
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 51671.332
$ws.Range("J24").Value = 65007
$ws.Range("L24").Value = 65007
$ws.Range("N24").Value = -65353

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1515.3
$ws.Range("I102").Value = 1819.375
$ws.Range("K102").Value = 1819.375
$ws.Range("M102").Value = -197.375

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7132.6665
$ws.Range("I126").Value = 6945.5
$ws.Range("K126").Value = 20836.5
$ws.Range("M126").Value = -18366.5

# GSM row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 104395.2
$ws.Range("J134").Value = 104395.2
$ws.Range("L134").Value = 313185.6
$ws.Range("N134").Value = -318255.6

# LTW row 45
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 7715
$ws.Range("I45").Value = 7715
$ws.Range("K45").Value = 7715
$ws.Range("M45").Value = -7308

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4656.9287
$ws.Range("I61").Value = 2171.2856
$ws.Range("K61").Value = 2171.2856
$ws.Range("M61").Value = -1969.2856

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5137.385
$ws.Range("I68").Value = 3169.4285
$ws.Range("J68").Value = 7433.3335
$ws.Range("K68").Value = 3169.4285
$ws.Range("L68").Value = 7433.3335
$ws.Range("M68").Value = -2420.4285
$ws.Range("N68").Value = -8931.333500000001

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5137.385
$ws.Range("I71").Value = 3169.4285
$ws.Range("J71").Value = 7433.3335
$ws.Range("K71").Value = 15847.1425
$ws.Range("L71").Value = 37166.6675
$ws.Range("M71").Value = -12103.1425
$ws.Range("N71").Value = -44654.6675

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4656.9287
$ws.Range("I113").Value = 2171.2856
$ws.Range("K113").Value = 2171.2856
$ws.Range("M113").Value = -1.285600000000159

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3333
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2999.5
$ws.Range("I132").Value = 4998
$ws.Range("J132").Value = 2333.3333
$ws.Range("K132").Value = 14994
$ws.Range("L132").Value = 6999.999899999999
$ws.Range("M132").Value = -12464
$ws.Range("N132").Value = -12059.9999
